$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("samples_retained")

# Update collected record counts for BAUM1 (row 4)
$ws.Range("C4").Value = 228
$ws.Range("D4").Value = 891

# Update the notes text for BAUM1 to mention dropped records
$ws.Range("H4").Value = "contempt, surprise, unsure, and boredom mapped to negative; labels determined by interrater consensus; some of the mp4s might not have audio!; interest mapped to positive; 61 records missing labels dropped"

# Update the selected cell in the sheet view
$ws.Range("C5").Select()

$excel.Calculate()
